# Updated relative_weights to v3
# Rearrange the weights: insert a new weighted column (F) - pushing the
# existing "Total weight" / "Weekly Goal" / "(Weight/Person)/week" columns
# one slot to the right - move "chanelling" out of column B into the new
# column F together with a new "if/then/else/fi" row, and retitle the old
# "if/then/else/fi" entry in column E to "redirection".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before the old "Total weight" column (F). Excel
#    clones the formatting of the column to its left (E) onto the new
#    column, and shifts F->G, G->H, H->I, fixing up the row-6 formulas.
$ws.Columns("F").Insert()

# 2. Fill in the new column F (weight value, header words, blank marker,
#    weekly count) - formatting was already inherited from column E above.
$ws.Range("F1").Value = 11
$ws.Range("F2").Value = "if/then/else/fi"
$ws.Range("F3").Value = "chanelling"
$ws.Range("F4").Value = "-"
$ws.Range("F5").Value = "-"
$ws.Range("F6").Value = 22

# 3. Update the rearranged labels in the existing columns.
$ws.Range("B2").Value = "-"
$ws.Range("E3").Value = "redirection"

# 4. Update the weekly counts in row 6 for the other columns.
$ws.Range("B6").Value = 0
$ws.Range("E6").Value = 21

# 5. Update the total-weight formula so it also includes the new column F.
$ws.Range("G6").Formula = "=(E6+D6+C6+B6+A6+F6)+9"

# 6. The right-most header cell ("Weekly Goal", now H5) loses its left
#    border now that it is no longer the first cell of that header group.
$ws.Range("H5").Borders.Item(7).LineStyle = -4142

# 7. Restore the active selection like the captured session (cell G7).
$ws.Range("G7").Select()
